$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 2675.1
$ws.Range("I86").Value = 2562.75
$ws.Range("J86").Value = 3124.5
$ws.Range("K86").Value = 2562.75
$ws.Range("L86").Value = 3124.5
$ws.Range("M86").Value = -1439.75
$ws.Range("N86").Value = -5370.5
$ws.Range("H89").Value = 2675.1
$ws.Range("I89").Value = 2562.75
$ws.Range("J89").Value = 3124.5
$ws.Range("K89").Value = 12813.75
$ws.Range("L89").Value = 15622.5
$ws.Range("M89").Value = -7197.75
$ws.Range("N89").Value = -26854.5
$ws.Range("H92").Value = 2686.1428
$ws.Range("I92").Value = 2826
$ws.Range("J92").Value = 2499.6667
$ws.Range("K92").Value = 2826
$ws.Range("L92").Value = 2499.6667
$ws.Range("M92").Value = -1578
$ws.Range("N92").Value = -4995.6667

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 9500
$ws.Range("J61").Value = 0
$ws.Range("L61").Value = 0
$ws.Range("N61").ClearContents()
$ws.Range("H74").Value = 307494.66
$ws.Range("I74").Value = 834860.5
$ws.Range("J74").Value = 6142.7617
$ws.Range("K74").Value = 834860.5
$ws.Range("L74").Value = 6142.7617
$ws.Range("M74").Value = -833986.5
$ws.Range("N74").Value = -7890.7617
$ws.Range("H77").Value = 307494.66
$ws.Range("I77").Value = 834860.5
$ws.Range("J77").Value = 6142.7617
$ws.Range("K77").Value = 4174302.5
$ws.Range("L77").Value = 30713.8085
$ws.Range("M77").Value = -4169934.5
$ws.Range("N77").Value = -39449.8085
$ws.Range("H122").Value = 1828.7142
$ws.Range("I122").Value = 1828.7142
$ws.Range("K122").Value = 5486.142599999999
$ws.Range("M122").Value = -3036.142599999999
$ws.Range("H136").Value = 9500
$ws.Range("J136").Value = 0
$ws.Range("L136").Value = 0
$ws.Range("N136").ClearContents()

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1185.8387
$ws.Range("I94").Value = 1175.963
$ws.Range("K94").Value = 1175.963
$ws.Range("M94").Value = -724.963

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 1686.1154
$ws.Range("I107").Value = 1447.8889
$ws.Range("K107").Value = 1447.8889
$ws.Range("M107").Value = 472.1111000000001

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H6").Value = 569.8
$ws.Range("I6").Value = 80
$ws.Range("J6").Value = 1304.5
$ws.Range("K6").Value = 80
$ws.Range("L6").Value = 1304.5
$ws.Range("M6").Value = 33
$ws.Range("N6").Value = -1530.5
$ws.Range("H7").Value = 914027.9399999999
$ws.Range("I7").Value = 2750
$ws.Range("J7").Value = 1116534.1
$ws.Range("K7").Value = 2750
$ws.Range("L7").Value = 1116534.1
$ws.Range("M7").Value = -2638
$ws.Range("N7").Value = -1116758.1
$ws.Range("H8").Value = 914027.9399999999
$ws.Range("I8").Value = 2750
$ws.Range("J8").Value = 1116534.1
$ws.Range("K8").Value = 2750
$ws.Range("L8").Value = 1116534.1
$ws.Range("M8").Value = -2611
$ws.Range("N8").Value = -1116812.1
$ws.Range("H12").Value = 114357.86
$ws.Range("I12").Value = 200048.75
$ws.Range("K12").Value = 200048.75
$ws.Range("M12").Value = -199908.75
$ws.Range("H14").Value = 352402
$ws.Range("I14").Value = 525298.7
$ws.Range("J14").Value = 6608.6
$ws.Range("K14").Value = 525298.7
$ws.Range("L14").Value = 6608.6
$ws.Range("M14").Value = -525130.7
$ws.Range("N14").Value = -6944.6
$ws.Range("H16").Value = 569.8
$ws.Range("I16").Value = 80
$ws.Range("J16").Value = 1304.5
$ws.Range("K16").Value = 80
$ws.Range("L16").Value = 1304.5
$ws.Range("M16").Value = 170
$ws.Range("N16").Value = -1804.5
$ws.Range("H126").Value = 2384.5151
$ws.Range("I126").Value = 2145.762
$ws.Range("J126").Value = 2802.3333
$ws.Range("K126").Value = 6437.286
$ws.Range("L126").Value = 8406.999899999999
$ws.Range("M126").Value = -3967.286
$ws.Range("N126").Value = -13346.9999
$ws.Range("H132").Value = 6204
$ws.Range("I132").Value = 3780
$ws.Range("J132").Value = 10809.6
$ws.Range("K132").Value = 11340
$ws.Range("L132").Value = 32428.8
$ws.Range("M132").Value = -8810
$ws.Range("N132").Value = -37488.8

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 92327.37
$ws.Range("I7").Value = 112522.336
$ws.Range("J7").Value = 1450
$ws.Range("K7").Value = 112522.336
$ws.Range("L7").Value = 1450
$ws.Range("M7").Value = -112410.336
$ws.Range("N7").Value = -1674
$ws.Range("H40").Value = 4717.1816
$ws.Range("I40").Value = 4991.1
$ws.Range("K40").Value = 4991.1
$ws.Range("M40").Value = -4855.1
$ws.Range("H46").Value = 6438.4116
$ws.Range("I46").Value = 2759.8
$ws.Range("K46").Value = 2759.8
$ws.Range("M46").Value = -2571.8
$ws.Range("H55").Value = 776.5294
$ws.Range("I55").Value = 496.3846
$ws.Range("J55").Value = 1687
$ws.Range("K55").Value = 496.3846
$ws.Range("L55").Value = 1687
$ws.Range("M55").Value = -323.3846
$ws.Range("N55").Value = -2033
$ws.Range("H61").Value = 7270.3887
$ws.Range("I61").Value = 7252.1333
$ws.Range("K61").Value = 7252.1333
$ws.Range("M61").Value = -7050.1333
$ws.Range("H93").Value = 1351.2
$ws.Range("I93").Value = 1345.7778
$ws.Range("K93").Value = 1345.7778
$ws.Range("M93").Value = -97.77780000000007
$ws.Range("H113").Value = 7270.3887
$ws.Range("I113").Value = 7252.1333
$ws.Range("K113").Value = 7252.1333
$ws.Range("M113").Value = -5082.1333
$ws.Range("H122").Value = 55563010
$ws.Range("I122").Value = 62506508
$ws.Range("K122").Value = 187519524
$ws.Range("M122").Value = -187517074
$ws.Range("H126").Value = 92327.37
$ws.Range("I126").Value = 112522.336
$ws.Range("J126").Value = 1450
$ws.Range("K126").Value = 337567.008
$ws.Range("L126").Value = 4350
$ws.Range("M126").Value = -335097.008
$ws.Range("N126").Value = -9290
$ws.Range("H133").Value = 69241.75
$ws.Range("J133").Value = 69241.75
$ws.Range("L133").Value = 69241.75
$ws.Range("N133").Value = -74301.75

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 6251.7646
$ws.Range("I81").Value = 1565
$ws.Range("K81").Value = 3130
$ws.Range("M81").Value = -2069
$ws.Range("H84").Value = 6251.7646
$ws.Range("I84").Value = 1565
$ws.Range("K84").Value = 15650
$ws.Range("M84").Value = -10346
$ws.Range("H113").Value = 731.17645
$ws.Range("I113").Value = 462.75
$ws.Range("K113").Value = 1388.25
$ws.Range("M113").Value = 781.75
$ws.Range("H122").Value = 3478.8057
$ws.Range("I122").Value = 3645.1333
$ws.Range("J122").Value = 2647.1667
$ws.Range("K122").Value = 10935.3999
$ws.Range("L122").Value = 7941.500100000001
$ws.Range("M122").Value = -8485.3999
$ws.Range("N122").Value = -12841.5001
$ws.Range("H126").Value = 5934.385
$ws.Range("I126").Value = 5377.091
$ws.Range("K126").Value = 16131.273
$ws.Range("M126").Value = -13661.273

Write-Output "All edits applied"